$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 2406.633596400077
$ws.Range("E2").Value = 290934.8942174729
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 148652.5872276
$ws.Range("L2").Value = 509125.9821312752
$ws.Range("M2").Value = 112470.9127927
$ws.Range("N2").Value = 72239.98700165171
$ws.Range("O2").Value = 68708.80120585459

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 53116.21837418353
$ws.Range("E2").Value = 269898.1793720496
$ws.Range("I2").Value = 222508.8162963558
$ws.Range("L2").Value = 224190.3804794976
$ws.Range("M2").Value = 105604.6794510125
$ws.Range("N2").Value = 35083.91479049736
$ws.Range("O2").Value = 25042.43351431981

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 12018.77283394627
$ws.Range("B2").Value = 20007.88192133839
$ws.Range("E2").Value = 138810.1661255918
$ws.Range("I2").Value = 161325.8753704437
$ws.Range("M2").Value = 64556.96661513005
$ws.Range("N2").Value = 45628.2541226642
$ws.Range("O2").Value = 58504.81512360305

# Sheet "2040" (sheet4.xml)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 264.8241227336373

# Sheet "2045" (sheet5.xml)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 41492.41011352674
